# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
# Both sheets contain the same event data and both received identical
# updates in the source commit.

$wb = $excel.ActiveWorkbook

# Map of row -> new value for column F (same updates on both sheets).
$updates = @{
    2  = 1097
    3  = 803
    6  = 1106
    8  = 2035
    9  = 7562
    11 = 413
    13 = 129
    14 = 399
    15 = 150
    16 = 7108
    17 = 299
    18 = 1333
    20 = 120
    22 = 141
    23 = 298
    24 = 136
    27 = 105
    34 = 76
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}

$wb.Save()
